# quotations/management/commands/oblupricelist.xlsx
# Split the combined "Tax Rate" column (F, stored as a 0..1 fraction) into
# two columns: F keeps a fraction-styled column now labelled "Tax Rate2",
# and a new column H is added, labelled "Tax Rate", holding the same rate
# expressed as a whole-number percentage (e.g. 0.18 -> 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the (renamed) existing column F.
$ws.Range("F1").Value = "Tax Rate2"

# Header for the new column H - reuses the original "Tax Rate" text.
$ws.Range("H1").Value = "Tax Rate"

# Last data row in the sheet.
$lastRow = 184

for ($r = 2; $r -le $lastRow; $r++) {
    $rate = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 8).Value = $rate * 100
}

# Move the active selection to H1, matching the saved workbook state.
$ws.Range("H1").Select() | Out-Null
